$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("n87", "n87_IMG_3088.jpeg"),
    @("n88", "n88_IMG_3088HorFlip.jpeg"),
    @("n89", "n89_IMG_3088HorVertFlip.jpeg"),
    @("n90", "n90_IMG_3088VertFlip.jpeg"),
    @("n91", "n91_IMG_3089.jpeg"),
    @("n92", "n92_IMG_3089HorFlip.jpeg"),
    @("n93", "n93_IMG_3089HorVertFlip.jpeg"),
    @("n94", "n94_IMG_3089VertFlip.jpeg"),
    @("n95", "n95_IMG_3091.jpeg"),
    @("n96", "n96_IMG_3091HorFlip.jpeg"),
    @("n97", "n97_IMG_3091HorVertFlip.jpeg"),
    @("n98", "n98_IMG_3091VertFlip.jpeg")
)

$startRow = 88
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $subjectId = $rows[$i][0]
    $fileName = $rows[$i][1]

    $ws.Cells.Item($r, 1).Value = $subjectId
    $ws.Cells.Item($r, 2).Value = $fileName
    $ws.Cells.Item($r, 3).Value = "'True"
    $ws.Cells.Item($r, 4).Value = "no_meltpatch"
    $ws.Cells.Item($r, 5).Value = "negative"
}

$wb.Save()
